$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new "Price" values look like plain decimal numbers (e.g. "1.00", "0.610").
# Excel would normally coerce such strings into numeric cells, dropping the
# significant trailing/leading zeros that are part of the displayed text in the
# source data. Force those specific cells to Text format first so the literal
# string is preserved, then restore the default "Normal" style afterward so no
# extra styling is left behind on the cell.
function Set-TextValue($ref, $value) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '63.600.70'
$ws.Range('E2').Value = '  -2.71%  '

$ws.Range('D3').Value = '3.319.33'
$ws.Range('E3').Value = '  -4.42%  '

$ws.Range('E4').Value = '  -0.12%  '

Set-TextValue 'D5' '548.21'
$ws.Range('E5').Value = '  -1.11%  '

Set-TextValue 'D6' '171.95'
$ws.Range('E6').Value = '  -4.18%  '

$ws.Range('E7').Value = '  -3.99%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D9' '0.610'
$ws.Range('E9').Value = '  -3.97%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D10' '0.151'
$ws.Range('E10').Value = '  -0.73%  '

$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D11' '53.17'
$ws.Range('E11').Value = '  -1.92%  '

$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D12' '0.0000264'
$ws.Range('E12').Value = '  -2.64%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '8.84'
$ws.Range('E13').Value = '  -4.34%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.867.11'
$ws.Range('E14').Value = '  -4.19%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D15' '18.13'
$ws.Range('E15').Value = '  -3.10%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.346.62'
$ws.Range('E16').Value = '  -3.77%  '

$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D17' '0.117'
$ws.Range('E17').Value = '  -3.70%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D18' '11.63'
$ws.Range('E18').Value = '  -3.21%  '

$ws.Range('D19').Value = '63.505.78'
$ws.Range('E19').Value = '  -3.19%  '

$ws.Range('B20').Value = 'Polygon'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D20' '0.973'
$ws.Range('E20').Value = '  -1.69%  '

$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D21' '410.53'
$ws.Range('E21').Value = '  -1.49%  '

$ws.Range('B22').Value = 'Toncoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D22' '4.39'
$ws.Range('E22').Value = '  +2.52%  '

Set-TextValue 'D23' '4.03'
$ws.Range('E23').Value = '  -0.97%  '

$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D24' '13.67'
$ws.Range('E24').Value = '  +6.95%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D25' '82.86'
$ws.Range('E25').Value = '  -3.61%  '

$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D26' '10.50'
$ws.Range('E26').Value = '  -3.52%  '

$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D27' '2.71'
$ws.Range('E27').Value = '  -5.58%  '

$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D28' '8.58'
$ws.Range('E28').Value = '  -5.76%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D29' '28.96'
$ws.Range('E29').Value = '  -4.62%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D30' '6.35'
$ws.Range('E30').Value = '  -3.49%  '

$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D31' '576.59'
$ws.Range('E31').Value = '  -5.80%  '

Set-TextValue 'D32' '11.31'
$ws.Range('E32').Value = '  -3.81%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.106'
$ws.Range('E33').Value = '  -3.86%  '

$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D34' '57.65'
$ws.Range('E34').Value = '  -2.51%  '

$ws.Range('B35').Value = 'Dai'
$ws.Range('C35').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D35' '1.00'
$ws.Range('E35').Value = '  -0.27%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D36' '0.147'
$ws.Range('E36').Value = '  +0.68%  '

$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D37' '34.92'
$ws.Range('E37').Value = '  -7.02%  '

$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D38' '3.39'
$ws.Range('E38').Value = '  +2.86%  '

$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '0.0₃0733'
$ws.Range('E39').Value = '  -7.24%  '

$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D40' '0.364'
$ws.Range('E40').Value = '  -4.57%  '

$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.114.90'
$ws.Range('E41').Value = '  -7.14%  '

$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D42' '1.00'
$ws.Range('E42').Value = '  -0.21%  '

$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D43' '2.77'
$ws.Range('E43').Value = '  -2.17%  '

$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 'D44' '3.23'
$ws.Range('E44').Value = '  -1.65%  '

$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D45' '0.0399'
$ws.Range('E45').Value = '  -3.85%  '

$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D46' '2.40'
$ws.Range('E46').Value = '  -5.78%  '

$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D47' '2.60'
$ws.Range('E47').Value = '  -4.44%  '

$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D48' '0.128'
$ws.Range('E48').Value = '  -3.87%  '

$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D49' '132.53'
$ws.Range('E49').Value = '  -3.78%  '

$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D50' '7.99'
$ws.Range('E50').Value = '  -5.29%  '

$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue 'D51' '0.000226'
$ws.Range('E51').Value = '  +4.89%  '
